$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.525.53'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.39%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.215.99'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.88%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '607.51'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.06%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '158.14'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.56%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.215.18'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.09%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.550'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.49%  '
$ws.Range('E10').Value = '  +0.58%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.68'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -4.20%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.502'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -2.78%  '
$ws.Range('E13').Value = '  +0.52%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '38.66'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.01%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.748.53'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.93%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.626.35'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.45%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '7.36'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.07%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.223.06'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.94%  '
$ws.Range('E19').Value = '  +1.15%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '506.57'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.51%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.14'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.38%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.732'
$ws.Range('D22').Style = 'Normal'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.00'
$ws.Range('D23').Style = 'Normal'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '14.57'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.68%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '84.83'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E26').Value = '  +0.11%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.99'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.27%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.10'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.44%  '
$ws.Range('E29').Value = '  +0.43%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.122'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +36.12%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.93'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.53%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.00'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.59%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '28.09'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.64%  '
$ws.Range('E34').Value = '  +0.14%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.18'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.21%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.47'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.58%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '501.30'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.59%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '55.35'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.97%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0₃0769'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +13.85%  '
$ws.Range('E40').Value = '  +6.03%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0420'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.67%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.04'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +6.00%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.71'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.93%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.297'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.43%  '
$ws.Range('E45').Value = '  -0.02%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.910.89'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.27%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '28.14'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.90%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.41'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.20%  '
$ws.Range('E50').Value = '  -1.03%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '122.47'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.74%  '
